# Fruta / hortaliza, semanal
# Insert a new weekly record row before the current row 25, shifting the
# existing rows 25-28 down to 26-29 (and extending the used range to R29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25; existing rows 25-28 shift down to 26-29.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new weekly record.
$ws.Cells.Item(25, 1).Value  = 1
$ws.Cells.Item(25, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(25, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(25, 4).Value  = 44776
$ws.Cells.Item(25, 5).Value  = 15
$ws.Cells.Item(25, 6).Value  = 100112043
$ws.Cells.Item(25, 7).Value  = "Pepino dulce"
$ws.Cells.Item(25, 8).Value  = "Cultivar IV Región"
$ws.Cells.Item(25, 9).Value  = "Primera"
$ws.Cells.Item(25, 10).Value = 200
$ws.Cells.Item(25, 11).Value = 17000
$ws.Cells.Item(25, 12).Value = 18000
$ws.Cells.Item(25, 13).Value = 17500
$ws.Cells.Item(25, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(25, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(25, 16).Value = 972
$ws.Cells.Item(25, 17).Value = 18
$ws.Cells.Item(25, 18).Value = "Hortaliza"

# Make sure the D25 cell keeps the date number format used throughout the
# "Fecha" column, matching the format of the row below it.
$ws.Cells.Item(25, 4).NumberFormat = $ws.Cells.Item(26, 4).NumberFormat
